$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 115, shifting existing rows 115-207 down to 116-208
$ws.Rows(115).Insert()

# Populate the newly inserted row 115 with the new data record
$ws.Cells.Item(115,1).Value = 7
$ws.Cells.Item(115,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(115,3).Value = "Ñuble"
$ws.Cells.Item(115,4).Value = 44566
$ws.Cells.Item(115,5).Value = 16
$ws.Cells.Item(115,6).Value = 100112002
$ws.Cells.Item(115,7).Value = "Pimiento"
$ws.Cells.Item(115,8).Value = "Cuatro cascos verde"
$ws.Cells.Item(115,9).Value = "Primera"
$ws.Cells.Item(115,10).Value = 120
$ws.Cells.Item(115,11).Value = 12000
$ws.Cells.Item(115,12).Value = 13000
$ws.Cells.Item(115,13).Value = 12500
$ws.Cells.Item(115,14).Value = "$/caja 15 kilos"
$ws.Cells.Item(115,15).Value = "Región del Maule"
$ws.Cells.Item(115,16).Value = 833
$ws.Cells.Item(115,17).Value = 15
$ws.Cells.Item(115,18).Value = "Hortaliza"
